# Auto-generated files on 2025-08-27
# Update the "同花顺" (column C) hot-stock list values for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value  = "领益智造"
$ws.Range("C6").Value  = "寒武纪"
$ws.Range("C11").Value = "吉视传媒"
$ws.Range("C12").Value = "大元泵业"
$ws.Range("C13").Value = "鸿博股份"
$ws.Range("C14").Value = "岩山科技"
$ws.Range("C15").Value = "工业富联"
$ws.Range("C16").Value = "合力泰"
$ws.Range("C17").Value = "新易盛"
$ws.Range("C18").Value = "大位科技"
$ws.Range("C20").Value = "东信和平"
$ws.Range("C21").Value = "启明信息"
